$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "37.162.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.058.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  +0.18%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "250.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.670"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  +1.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "61.31"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  +13.07%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -0.02%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  +2.71%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  +0.90%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  +2.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "16.27"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  +8.34%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  +0.19%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  +1.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "5.71"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  +9.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "2.063.54"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "18.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  +29.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "37.166.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "75.69"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  +4.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "0.0{0}0913" -f $sub3
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -3.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "5.47"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  +2.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "239.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  +1.22%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  -0.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  +13.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "169.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  -0.23%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  +4.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "20.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  +0.66%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  +2.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  +10.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.84"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  +6.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.0623"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  +0.63%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  +4.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.0888"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  +4.15%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  -0.07%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  +0.53%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  -2.04%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  +5.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.35"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  +1.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "5.36"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  +32.93%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  +14.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "18.19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  +1.21%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  +1.39%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  +0.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "98.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  +2.03%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  +3.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.303.38"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  +0.83%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  -1.11%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  +1.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.249.21"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  +0.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "3.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -10.75%  "
